# A new September transaction log entry was recorded in the "2024" sheet.
# It lands at the top of the September_Details/September_Date block (row 35
# of the "Others" group), pushing every row below it down by one - which in
# turn pushes the August block's remaining row into the Broadband group
# label row, and that label itself down into a brand-new row 108.
#
# Excel's native Rows().Insert() reproduces exactly that "shift everything
# down" behaviour (and grows the sheet's used range from A1:Y107 to
# A1:Y108), so we just insert a blank row at 35 and then fill in the new
# entry's description/timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new blank row at 35; rows 35-107 shift down to 36-108.
$ws.Rows("35:35").Insert()

# Fill in the new, most-recent September entry.
$ws.Range("R35").Value = "corporate internet share"
$ws.Range("S35").Value = "2024-09-09 11:35:34"
